$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.01.2022 18:15"

# Row 7 (MOL Olomoucka) price update: B7/C7 swap-like update, D7 delta, E7 timestamp text
$ws.Range("B7").Value = 36.9
$ws.Range("C7").Value = 36.7

# D7 must stay a literal text cell holding "+0.2" (not get auto-parsed into
# the number 0.2, and without minting a new number-format style on the
# cell). Writing the text through a scratch cell's formula result and
# pasting-special just the *values* back onto D7 preserves both the
# literal text and the cell's original (default) style.
$ws.Range("Z1").Formula = '="+0.2"'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null

$ws.Range("E7").Value = "2022-01-21 18:15:22"
